$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 247 (pushes the existing rows 247-308 down to 248-309)
$ws.Range("A247").EntireRow.Insert()

# Populate the newly inserted row with the new record
$ws.Range("A247").Value = 5
$ws.Range("B247").Value = "Macroferia Regional de Talca"
$ws.Range("C247").Value = "Maule"
$ws.Range("D247").Value = 44722
$ws.Range("E247").Value = 7
$ws.Range("F247").Value = 100112003
$ws.Range("G247").Value = "Ajo"
$ws.Range("H247").Value = "Chino"
$ws.Range("I247").Value = "Primera"
$ws.Range("J247").Value = 300
$ws.Range("K247").Value = 18000
$ws.Range("L247").Value = 18000
$ws.Range("M247").Value = 18000
$ws.Range("N247").Value = "$/caja 10 kilos"
$ws.Range("O247").Value = "China"
$ws.Range("P247").Value = 1800
$ws.Range("Q247").Value = 10
$ws.Range("R247").Value = "Hortaliza"
